$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "324.12"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.74%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.23"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.25%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.697"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "7.89%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08011"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.10%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.491"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.65%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.618"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.03%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.973"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.43%"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.84%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9283"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.85%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1240"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-6.27%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1974"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.16%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.707"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "24.74%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09158"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.56%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03608"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.94%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.1048"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "9.58%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001301"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.84%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006164"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.17%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.76%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3471"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.46%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.64%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2412"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-5.88%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04411"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.00%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001262"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.40%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004664"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.33%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-3.42%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02497"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.11%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05355"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.74%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007476"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.35%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009615"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.16%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1404"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.58%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002117"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.55%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01031"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.67%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006736"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.29%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.06%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002969"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-11.22%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002291"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-4.64%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.06%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.06%"
